$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Move the "select a role" validation-list helper column from Y to X
#    (Copy preserves values/formatting exactly, including the leading
#    apostrophe text-qualifiers, without Excel's "smart quote" reinterpretation.)
$src = $ws.Range("Y1:Y16")
$dst = $ws.Range("X1:X16")
$src.Copy($dst)
$ws.Range("Y1:Y21").ClearContents()

# 2. Point the dropdown validation on D2:D21 at the new X column list
$ws.Range("D2:D21").Validation.Formula1 = "=`$X`$1:`$X`$16"

# 3. Everyone gets the "admin" role for this simulation
$ws.Range("D2:D21").Value = "admin"

# 4. Restore the view to show column A first, matching the saved workbook
$ws.Range("A1").Select()
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("M1").Select()
